$p = $ppt.ActivePresentation

# --- Slide 11: "Use Return Path TLV for STAMP from the probe query message" ---
# Remove bold, add blue (0070C0) solid fill color to the run's font.
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(3)
$tr11 = $shp11.TextFrame.TextRange
$para11 = $tr11.Paragraphs(6, 1)
$para11.Font.Bold = $false
$para11.Font.Color.RGB = 12611584

# --- Slide 13: remove the "Implementation exists" paragraph ---
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(2)
$tr13 = $shp13.TextFrame.TextRange
$para13 = $tr13.Paragraphs(2, 1)
$para13.Delete()
